$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This export was regenerated with refreshed TPM numbers. The old rows 2-4
# (sending cluster "ECs") are gone; the old rows 5-7 (sending cluster "FAPs")
# move up into rows 2-4 with updated specificity/weight figures, and the
# now-unused trailing rows are removed.
$ws.Rows("5:7").Delete()

# Row 2: FAPs -> ECs (Cxcl13/Cxcr5), refreshed TPM numbers
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Cxcr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.91523966666667
$ws.Range("H2").Value = 32.745719
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05950833333333333
$ws.Range("N2").Value = 0.178525
$ws.Range("O2").Value = 0.04031524136301987
$ws.Range("P2").Value = 0.04031524136301987
$ws.Range("Q2").Value = 0.6495477204972222
$ws.Range("R2").Value = 5.845929484475
$ws.Range("S2").Value = 0.04031524136301987
$ws.Range("T2").Value = 0.04031524136301987

# Row 3: FAPs -> FAPs (Cxcl13/Cxcr5), refreshed TPM numbers
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Cxcr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.91523966666667
$ws.Range("H3").Value = 32.745719
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.181787
$ws.Range("N3").Value = 3.545361
$ws.Range("O3").Value = 0.8006278360679875
$ws.Range("P3").Value = 0.8006278360679875
$ws.Range("Q3").Value = 12.899488339951
$ws.Range("R3").Value = 116.095395059559
$ws.Range("S3").Value = 0.8006278360679875
$ws.Range("T3").Value = 0.8006278360679875

# Row 4: FAPs -> MuSCs (Cxcl13/Cxcr5), refreshed TPM numbers
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Cxcr5"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.91523966666667
$ws.Range("H4").Value = 32.745719
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.23478
$ws.Range("N4").Value = 0.70434
$ws.Range("O4").Value = 0.1590569225689926
$ws.Range("P4").Value = 0.1590569225689926
$ws.Range("Q4").Value = 2.56267996894
$ws.Range("R4").Value = 23.06411972046
$ws.Range("S4").Value = 0.1590569225689926
$ws.Range("T4").Value = 0.1590569225689926
